# Updates the cryptos list (prices / 1h volume %, plus a Cosmos/Hedera row
# swap) to match the latest scrape, as produced by the GitHub Actions job.
#
# The Price column (D) holds text that sometimes looks like a plain number
# (e.g. "1.00", "586.02"). Assigning such a string straight to .Value would
# make Excel auto-convert it to a numeric cell (losing formatting / trailing
# zeros and introducing floating point noise). To keep these as text, each
# Price cell is temporarily switched to the "@" (Text) number format before
# the assignment, then its style is restored to "Normal" afterwards so the
# cell's style index stays exactly as it was before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.417.62'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -4.95%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.249.63'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -8.23%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.02'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.86'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -12.48%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.241.42'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -8.27%  '
$ws.Range('E9').Value = '  -10.87%  '
$ws.Range('E10').Value = '  -12.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.82'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.38%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.507'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -13.65%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.59'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -17.46%  '
$ws.Range('E14').Value = '  -11.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.768.87'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -8.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.491.83'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.87%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.251.55'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -8.37%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '545.47'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -10.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.26'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -14.04%  '
$ws.Range('E20').Value = '  -5.83%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.24'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -14.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.767'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -13.63%  '
$ws.Range('E23').Value = '  -12.98%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.66'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -13.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.52'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -14.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('E27').Value = '  -15.31%  '
$ws.Range('E28').Value = '  -10.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '29.47'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -13.30%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.13'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -17.72%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.73'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -10.74%  '
$ws.Range('E32').Value = '  -11.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '553.82'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -12.48%  '
$ws.Range('E34').Value = '  -19.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.77'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -16.16%  '
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '53.84'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0442'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.80%  '
$ws.Range('B39').Value = 'Cosmos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '9.24'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -15.00%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0851'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -15.15%  '
$ws.Range('E41').Value = '  -11.47%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.940.60'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -12.68%  '
$ws.Range('E43').Value = '  -25.49%  '
$ws.Range('E44').Value = '  -16.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0₃0586'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -21.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '26.42'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -18.07%  '
$ws.Range('E47').Value = '  -20.43%  '
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.13'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -17.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '125.58'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.68%  '
